# Microservice Exploration Project - Progress.docx edit script
# Implements: "Planned things for order-service"
#
# Strategy: locate each target paragraph (by its plain text, ignoring the
# trailing paragraph mark) and replace its content with a freshly built
# OOXML fragment reflecting the desired end state (new runs, proofErr
# spell-check wrappers, relocated bookmark / lastRenderedPageBreak, etc).
# Using Range.InsertXML() on the whole paragraph range cleanly swaps the
# paragraph contents without leaving stray empty paragraphs behind.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

function Get-ParaByText {
    param($doc, [string]$text)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $doc.Paragraphs($i)
        }
    }
    return $null
}

function Get-ParaRangeByText {
    param($doc, [string]$startText, [string]$endText)
    $p1 = Get-ParaByText $doc $startText
    $p2 = Get-ParaByText $doc $endText
    if ($p1 -eq $null -or $p2 -eq $null) {
        throw "Could not locate paragraph(s): [$startText] .. [$endText]"
    }
    return $doc.Range($p1.Range.Start, $p2.Range.End)
}

function Replace-Para {
    param($doc, [string]$text, [string]$xml)
    $p = Get-ParaByText $doc $text
    if ($p -eq $null) {
        throw "Could not locate paragraph: [$text]"
    }
    $p.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1. Sprint 6 / 7 / 8 restructuring
#    Before: "Sprint 6 - Add API gateway service" (bookmark _GoBack after
#            "API"), "<tab>Added", "Sprint 7 - "
#    After:  bookmark relocated to new "Sprint 8 - " paragraph; new
#            sentences added for sprints 7 and 8.
# ---------------------------------------------------------------------
$sprintRange = Get-ParaRangeByText $d "Sprint 6 - Add API gateway service" "Sprint 7 - "

$sprintXml = @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Sprint 6 - </w:t></w:r>
  <w:r><w:t xml:space="preserve">Add </w:t></w:r>
  <w:r><w:t>API</w:t></w:r>
  <w:r><w:t xml:space="preserve"> gateway service</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r><w:tab/><w:t>Added</w:t></w:r>
  <w:r><w:t xml:space="preserve"> service with spring cloud gateway</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Sprint 7 </w:t></w:r>
  <w:r><w:t>&#8211;</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>Add monitoring services</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Sprint 8 - </w:t></w:r>
  <w:bookmarkStart w:id="1" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="1"/>
</w:p>
"@
$sprintRange.InsertXML($sprintXml)

# ---------------------------------------------------------------------
# 2. "WebFlux" hyperlink gains spell-check proofErr wrapping
# ---------------------------------------------------------------------
Replace-Para $d "Also (later): Reactive stream - WebFlux" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Also </w:t></w:r>
  <w:r><w:t>(</w:t></w:r>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>later</w:t></w:r>
  <w:r><w:t xml:space="preserve">): Reactive stream - </w:t></w:r>
  <w:hyperlink r:id="rId5" w:history="1">
    <w:proofErr w:type="spellStart"/>
    <w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>WebFlux</w:t></w:r>
    <w:proofErr w:type="spellEnd"/>
  </w:hyperlink>
</w:p>
"@

# ---------------------------------------------------------------------
# 3. "(try gRPC LATER)" split into 3 runs with proofErr around gRPC;
#    lastRenderedPageBreak moves from "UI: Angular" paragraph to the
#    "Database:" paragraph that follows.
# ---------------------------------------------------------------------
Replace-Para $d "Service communication: REST APIs with HATEOAS links (try gRPC LATER)" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
  </w:pPr>
  <w:r><w:t>Service communication: REST APIs</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve">with </w:t></w:r>
  <w:r><w:t xml:space="preserve">HATEOAS </w:t></w:r>
  <w:r><w:t xml:space="preserve">links </w:t></w:r>
  <w:r><w:t xml:space="preserve">(try </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>gRPC</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> LATER)</w:t></w:r>
</w:p>
"@

Replace-Para $d "Database: In-memory H2 for now, integrate Postgres & NoSQL(MongoDB/DynamoDB)  later" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
  </w:pPr>
  <w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Database: </w:t></w:r>
  <w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>In-memory H2 for now</w:t></w:r>
  <w:r><w:t>, integrate Postgres</w:t></w:r>
  <w:r><w:t xml:space="preserve"> &amp; NoSQL(MongoDB/DynamoDB) </w:t></w:r>
  <w:r><w:t xml:space="preserve"> later</w:t></w:r>
</w:p>
"@

Replace-Para $d "UI: Angular 8, TypeScript, Bootstrap" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">UI: Angular </w:t></w:r>
  <w:r><w:t>8</w:t></w:r>
  <w:r><w:t>, TypeScript, Bootstrap</w:t></w:r>
</w:p>
"@

# ---------------------------------------------------------------------
# 4. "Micrometer" hyperlink gains proofErr wrapping
# ---------------------------------------------------------------------
Replace-Para $d "Monitoring: Micrometer + Prometheus" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Monitoring: </w:t></w:r>
  <w:hyperlink r:id="rId10" w:history="1">
    <w:proofErr w:type="spellStart"/>
    <w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:b/><w:bCs/></w:rPr><w:t>Micrometer</w:t></w:r>
    <w:proofErr w:type="spellEnd"/>
  </w:hyperlink>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> + </w:t></w:r>
  <w:hyperlink r:id="rId11" w:history="1">
    <w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:b/><w:bCs/></w:rPr><w:t>Prometheus</w:t></w:r>
  </w:hyperlink>
</w:p>
"@

# ---------------------------------------------------------------------
# 5. "Zipkin" hyperlink gains proofErr wrapping
# ---------------------------------------------------------------------
Replace-Para $d "Request Tracing:   Spring Cloud Sleuth + Zipkin" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Request </w:t></w:r>
  <w:hyperlink r:id="rId12" w:history="1">
    <w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Tracing</w:t></w:r>
  </w:hyperlink>
  <w:r><w:t xml:space="preserve">: </w:t></w:r>
  <w:r><w:t xml:space="preserve">  </w:t></w:r>
  <w:hyperlink r:id="rId13" w:tgtFrame="_blank" w:history="1">
    <w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:b/><w:bCs/></w:rPr><w:t>Spring Cloud Sleuth</w:t></w:r>
  </w:hyperlink>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> + </w:t></w:r>
  <w:hyperlink r:id="rId14" w:history="1">
    <w:proofErr w:type="spellStart"/>
    <w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:b/><w:bCs/></w:rPr><w:t>Zipkin</w:t></w:r>
    <w:proofErr w:type="spellEnd"/>
  </w:hyperlink>
</w:p>
"@

# ---------------------------------------------------------------------
# 6. "Spring Cloud LoadBalancer" hyperlink split, proofErr around
#    "LoadBalancer"
# ---------------------------------------------------------------------
Replace-Para $d "Load Balancing: Spring Cloud LoadBalancer" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Load Balancing: </w:t></w:r>
  <w:hyperlink r:id="rId17" w:history="1">
    <w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Spring Cloud </w:t></w:r>
    <w:proofErr w:type="spellStart"/>
    <w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:b/><w:bCs/></w:rPr><w:t>LoadBalancer</w:t></w:r>
    <w:proofErr w:type="spellEnd"/>
  </w:hyperlink>
</w:p>
"@

# ---------------------------------------------------------------------
# 7. "DevTools" gains proofErr wrapping
# ---------------------------------------------------------------------
Replace-Para $d "	JPA, Swagger2 (for API docs), DevTools, Actuator, " @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
  </w:pPr>
  <w:r><w:tab/><w:t xml:space="preserve">JPA, </w:t></w:r>
  <w:r><w:t>Swagger</w:t></w:r>
  <w:r><w:t>2</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve">(for API docs), </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>DevTools</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">, Actuator, </w:t></w:r>
</w:p>
"@

# ---------------------------------------------------------------------
# 8. "Unit: @DataJpa, @WebMvcTest, @SpringBootTest" - split & proofErr
# ---------------------------------------------------------------------
Replace-Para $d "Unit: @DataJpa, @WebMvcTest, @SpringBootTest" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r><w:t>Unit: @</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>DataJpa</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>, @</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Web</w:t></w:r>
  <w:r><w:t>MvcTest</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> @</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>SpringBootTest</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@

# ---------------------------------------------------------------------
# 9. "/CloudFoundry" split with proofErr around "CloudFoundry"
# ---------------------------------------------------------------------
Replace-Para $d "=> Maybe try on OpenStack/CloudFoundry as well at last" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:ind w:left="720" w:firstLine="720"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">=&gt; </w:t></w:r>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Maybe </w:t></w:r>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">try </w:t></w:r>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">on </w:t></w:r>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>OpenStack</w:t></w:r>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>/</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>CloudFoundry</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> as well</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>at last</w:t></w:r>
</w:p>
"@

# ---------------------------------------------------------------------
# 10. lastRenderedPageBreak moves from "Other Editor Tools:" paragraph
#     up to the "SCM:  GitHub" paragraph.
# ---------------------------------------------------------------------
Replace-Para $d "SCM:  GitHub" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
  </w:pPr>
  <w:r><w:lastRenderedPageBreak/><w:t>SCM:  GitHub</w:t></w:r>
</w:p>
"@

Replace-Para $d "Other Editor Tools: Git Bash, Eclipse, Visual Studio Code, Postman, SourceTree" @"
<w:p $wNs>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Other Editor Tools: </w:t></w:r>
  <w:r><w:t xml:space="preserve">Git Bash, </w:t></w:r>
  <w:r><w:t>Eclipse, Visual Studio Code, Postman</w:t></w:r>
  <w:r><w:t>, SourceTree</w:t></w:r>
</w:p>
"@
